$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CC-Helping")
$ws.Activate()

# Rename header cell A1 from "DL_Cond" to "CashOrCredit_Ferret"
$ws.Range("A1").Value = "CashOrCredit_Ferret"

# Replace the numeric 1/2 coding in column A with text labels "credit_card"/"cash"
$ws.Range("A2").Formula = '=IF(ISODD(B2),"credit_card", "cash")'
$ws.Range("A3:A66").Formula = '=IF(ISODD(B3),"credit_card", "cash")'
$ws.Range("A67:A95").Formula = '=IF(ISODD(B67),"credit_card", "cash")'

# Widen column A to fit the new text labels (column B keeps its default width)
$ws.Columns.Item(1).ColumnWidth = 20

# Scroll the frozen pane back up to the top and change the active selection
$w = $excel.ActiveWindow
$w.ScrollRow = 2
$w.ScrollColumn = 1
$ws.Range("B1:B1048576").Select()
